$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 277, pushing the existing rows 277-297 down to 278-298
# (carries formatting, e.g. the date style on column D, along with the rows).
$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new weekly record.
$ws.Range("A277").Value = 4
$ws.Range("B277").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C277").Value = "Los Lagos"
$ws.Range("D277").Value = 44769
$ws.Range("E277").Value = 10
$ws.Range("F277").Value = 100112043
$ws.Range("G277").Value = "Pepino ensalada"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 70
$ws.Range("K277").Value = 22000
$ws.Range("L277").Value = 22000
$ws.Range("M277").Value = 22000
$ws.Range("N277").Value = "$/caja 60 unidades"
$ws.Range("O277").Value = "Región de Arica y Parinacota"
$ws.Range("P277").Value = 367
$ws.Range("Q277").Value = 60
$ws.Range("R277").Value = "Hortaliza"
